$wb = $excel.ActiveWorkbook

# Sheets involved (1-based index, matching workbook.xml <sheets> order):
#  4 = "ЩО-0.1-81"  (currently the active/selected tab)
#  6 = "ЩО-0.3-83"  (selection needs to move to B17)
#  7 = "ЩО-1.1-84"  (new "Sun Time"/"outdoor" group row is added here; becomes active tab)

$ws4 = $wb.Worksheets.Item(4)
$ws6 = $wb.Worksheets.Item(6)
$ws7 = $wb.Worksheets.Item(7)

# --- Insert the new "Sun Time" / "outdoor" group as row 6 of ЩО-1.1-84, ---
# --- shifting the existing rows (old 6-21) down to 7-22.                ---
$ws7.Rows("6:6").Insert(-4121)

# Give the new row the same cell formatting as the data row right below it
# (which is the row that used to be row 6).
$ws7.Range("A7:E7").Copy()
$ws7.Range("A6:E6").PasteSpecial(-4122)

# Column A ("Helvar group number") holds the text "10" (not the number 10) for
# this row, just like the existing text cell on the ЩО-0.1-81 sheet. Copy it
# over as a value so it keeps its text type instead of being read as numeric.
$ws4.Range("C19").Copy()
$ws7.Range("A6").PasteSpecial(-4163)

$ws7.Range("B6").Value = "Sun Time"
$ws7.Range("C6").Value = "outdoor"

# Column D ("Dimming") needs the text "false" (not the boolean FALSE). Reuse a
# neighboring cell that already stores "false" as text and copy its value.
$ws7.Range("D8").Copy()
$ws7.Range("D6").PasteSpecial(-4163)

$ws7.Range("E6").Value = 0

# The "Dimming" column's list data-validation covered D6:D1048576 before the
# insert pushed it down to D7:D1048576; restore it back to start at row 6 so
# the new row is included too.
$ws7.Range("D7:D1048576").Validation.Delete()
$ws7.Range("D6:D1048576").Validation.Add(3, 1, 1, "диммирование")

# --- Update stored selections on the sheets that had their view touched. ---
$ws6.Range("B17").Select()

# --- Switch the active tab from ЩО-0.1-81 to ЩО-1.1-84. ---
$ws7.Activate()
$ws7.Range("C12").Select()
